# Refresh market-price-derived Leve profit columns (H:N) for specific leve rows
# across several job sheets, per the scheduled-runner data update.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
# M/N are only present when the corresponding NQ/HQ price data exists, so some
# rows need a cell cleared (no longer applicable) or newly set (now applicable).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Cells.Item(6, 8).Value = 357214.28
$ws.Cells.Item(6, 9).Value = 416745.4
$ws.Cells.Item(6, 10).Value = 27.5
$ws.Cells.Item(6, 11).Value = 1250236.2
$ws.Cells.Item(6, 12).Value = 82.5
$ws.Cells.Item(6, 13).Value = -1250124.2
$ws.Cells.Item(6, 14).Value = -306.5
# ALC row 18
$ws.Cells.Item(18, 8).Value = 2720.1667
$ws.Cells.Item(18, 9).Value = 2266.4
$ws.Cells.Item(18, 10).Value = 4989
$ws.Cells.Item(18, 11).Value = 2266.4
$ws.Cells.Item(18, 12).Value = 4989
$ws.Cells.Item(18, 13).Value = -1982.4
$ws.Cells.Item(18, 14).Value = -5557
# ALC row 39
$ws.Cells.Item(39, 8).Value = 174.81818
$ws.Cells.Item(39, 9).Value = 192
$ws.Cells.Item(39, 10).Value = 3
$ws.Cells.Item(39, 11).Value = 576
$ws.Cells.Item(39, 12).Value = 9
$ws.Cells.Item(39, 13).Value = -280
$ws.Cells.Item(39, 14).Value = -601
# ALC row 41
$ws.Cells.Item(41, 8).Value = 265.1111
$ws.Cells.Item(41, 9).Value = 255.28572
$ws.Cells.Item(41, 10).Value = 299.5
$ws.Cells.Item(41, 11).Value = 255.28572
$ws.Cells.Item(41, 12).Value = 299.5
$ws.Cells.Item(41, 13).Value = 184.71428
$ws.Cells.Item(41, 14).Value = -1179.5
# ALC row 43
$ws.Cells.Item(43, 8).Value = 4000
$ws.Cells.Item(43, 10).Value = 3000
$ws.Cells.Item(43, 12).Value = 3000
$ws.Cells.Item(43, 14).Value = -3138
# ALC row 55
$ws.Cells.Item(55, 8).Value = 4000
$ws.Cells.Item(55, 9).Value = 4000
$ws.Cells.Item(55, 11).Value = 4000
$ws.Cells.Item(55, 13).Value = -3786
# ALC row 96
$ws.Cells.Item(96, 8).Value = 7999
$ws.Cells.Item(96, 9).Value = 7999
$ws.Cells.Item(96, 11).Value = 23997
$ws.Cells.Item(96, 13).Value = -22624
# ALC row 138
$ws.Cells.Item(138, 8).Value = 2310.5715
$ws.Cells.Item(138, 9).Value = 925.7143
$ws.Cells.Item(138, 10).Value = 3695.4285
$ws.Cells.Item(138, 11).Value = 2777.1429
$ws.Cells.Item(138, 12).Value = 11086.2855
$ws.Cells.Item(138, 13).Value = 2362.8571
$ws.Cells.Item(138, 14).Value = -21366.2855
$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Cells.Item(20, 8).Value = 3268.375
$ws.Cells.Item(20, 9).Value = 2269.3333
$ws.Cells.Item(20, 10).Value = 3867.8
$ws.Cells.Item(20, 11).Value = 2269.3333
$ws.Cells.Item(20, 12).Value = 3867.8
$ws.Cells.Item(20, 13).Value = -2022.3333
$ws.Cells.Item(20, 14).Value = -4361.8
# BSM row 76
$ws.Cells.Item(76, 8).Value = 21998
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 21998
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 21998
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 14).Value = -22628
# BSM row 79
$ws.Cells.Item(79, 8).Value = 21998
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 21998
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 21998
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(79, 14).Value = -24182
# BSM row 86
$ws.Cells.Item(86, 8).Value = 885.8889
$ws.Cells.Item(86, 9).Value = 713.4
$ws.Cells.Item(86, 10).Value = 1101.5
$ws.Cells.Item(86, 11).Value = 713.4
$ws.Cells.Item(86, 12).Value = 1101.5
$ws.Cells.Item(86, 13).Value = 409.6
$ws.Cells.Item(86, 14).Value = -3347.5
# BSM row 89
$ws.Cells.Item(89, 8).Value = 885.8889
$ws.Cells.Item(89, 9).Value = 713.4
$ws.Cells.Item(89, 10).Value = 1101.5
$ws.Cells.Item(89, 11).Value = 3567
$ws.Cells.Item(89, 12).Value = 5507.5
$ws.Cells.Item(89, 13).Value = 2049
$ws.Cells.Item(89, 14).Value = -16739.5
# BSM row 94
$ws.Cells.Item(94, 8).Value = 518.5
$ws.Cells.Item(94, 9).Value = 518.5
$ws.Cells.Item(94, 11).Value = 518.5
$ws.Cells.Item(94, 13).Value = -67.5
# BSM row 132
$ws.Cells.Item(132, 8).Value = 124000
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).ClearContents()
$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Cells.Item(22, 8).Value = 194.66667
$ws.Cells.Item(22, 9).Value = 194.66667
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 194.66667
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 155.33333
$ws.Cells.Item(22, 14).ClearContents()
# CRP row 68
$ws.Cells.Item(68, 8).Value = 47998
$ws.Cells.Item(68, 10).Value = 47998
$ws.Cells.Item(68, 12).Value = 47998
$ws.Cells.Item(68, 14).Value = -49496
# CRP row 71
$ws.Cells.Item(71, 8).Value = 47998
$ws.Cells.Item(71, 10).Value = 47998
$ws.Cells.Item(71, 12).Value = 143994
$ws.Cells.Item(71, 14).Value = -151482
# CRP row 86
$ws.Cells.Item(86, 8).Value = 9585.166999999999
$ws.Cells.Item(86, 9).Value = 9001
$ws.Cells.Item(86, 11).Value = 9001
$ws.Cells.Item(86, 13).Value = -7878
# CRP row 89
$ws.Cells.Item(89, 8).Value = 9585.166999999999
$ws.Cells.Item(89, 9).Value = 9001
$ws.Cells.Item(89, 11).Value = 45005
$ws.Cells.Item(89, 13).Value = -39389
$ws = $wb.Worksheets.Item("CUL")
# CUL row 11
$ws.Cells.Item(11, 8).Value = 3911.111
$ws.Cells.Item(11, 10).Value = 3911.111
$ws.Cells.Item(11, 12).Value = 11733.333
$ws.Cells.Item(11, 14).Value = -12013.333
# CUL row 55
$ws.Cells.Item(55, 8).Value = 1893.125
$ws.Cells.Item(55, 10).Value = 1893.125
$ws.Cells.Item(55, 12).Value = 5679.375
$ws.Cells.Item(55, 14).Value = -6033.375
# CUL row 80
$ws.Cells.Item(80, 9).Value = 2000
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 6000
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -5064
$ws.Cells.Item(80, 14).ClearContents()
# CUL row 83
$ws.Cells.Item(83, 9).Value = 2000
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 18000
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -13320
$ws.Cells.Item(83, 14).ClearContents()
# CUL row 109
$ws.Cells.Item(109, 8).Value = 4005.5
$ws.Cells.Item(109, 9).Value = 4005.5
$ws.Cells.Item(109, 11).Value = 12016.5
$ws.Cells.Item(109, 13).Value = -10976.5
# CUL row 128
$ws.Cells.Item(128, 8).Value = 295000
$ws.Cells.Item(128, 9).Value = 295000
$ws.Cells.Item(128, 11).Value = 885000
$ws.Cells.Item(128, 13).Value = -880020
# CUL row 129
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 13).ClearContents()
$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Cells.Item(80, 8).Value = 3588.75
$ws.Cells.Item(80, 9).Value = 2368.6667
$ws.Cells.Item(80, 11).Value = 2368.6667
$ws.Cells.Item(80, 13).Value = -1370.6667
# GSM row 83
$ws.Cells.Item(83, 8).Value = 3588.75
$ws.Cells.Item(83, 9).Value = 2368.6667
$ws.Cells.Item(83, 11).Value = 11843.3335
$ws.Cells.Item(83, 13).Value = -6851.333500000001
# GSM row 102
$ws.Cells.Item(102, 8).Value = 3745.2
$ws.Cells.Item(102, 9).Value = 3745.2
$ws.Cells.Item(102, 11).Value = 3745.2
$ws.Cells.Item(102, 13).Value = -2123.2
$ws = $wb.Worksheets.Item("LTW")
# LTW row 46
$ws.Cells.Item(46, 8).Value = 3620
$ws.Cells.Item(46, 9).Value = 4025
$ws.Cells.Item(46, 10).Value = 2000
$ws.Cells.Item(46, 11).Value = 4025
$ws.Cells.Item(46, 12).Value = 2000
$ws.Cells.Item(46, 13).Value = -3837
$ws.Cells.Item(46, 14).Value = -2376
